$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 39, shifting existing rows 39:87 down to 40:88.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new record's data.
$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(39, 3).Value = "Coquimbo"
$ws.Cells.Item(39, 4).Value = "10/14/2021"
$ws.Cells.Item(39, 5).Value = 4
$ws.Cells.Item(39, 6).Value = 100112040
$ws.Cells.Item(39, 7).Value = "Cilantro"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 3160
$ws.Cells.Item(39, 11).Value = 1500
$ws.Cells.Item(39, 12).Value = 2000
$ws.Cells.Item(39, 13).Value = 1750
$ws.Cells.Item(39, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(39, 16).Value = 1167
$ws.Cells.Item(39, 17).Value = 1.5
$ws.Cells.Item(39, 18).Value = "Hortaliza"
